# Scheduled market-data refresh: updates computed leve-profit figures
# (currentAveragePrice / NQ / HQ, LevePrice, LeveProfit columns) per sheet
# using latest market board snapshot values.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 454.35
$ws.Range("I33").Value = 443.86667
$ws.Range("K33").Value = 443.86667
$ws.Range("M33").Value = -214.86667
$ws.Range("H43").Value = 9278068
$ws.Range("I43").Value = 100001
$ws.Range("J43").Value = 11113681
$ws.Range("K43").Value = 100001
$ws.Range("L43").Value = 11113681
$ws.Range("M43").Value = -99932
$ws.Range("N43").Value = -11113819
$ws.Range("H129").Value = 837.7317
$ws.Range("I129").Value = 563.75
$ws.Range("K129").Value = 1691.25
$ws.Range("M129").Value = 3308.75
$ws.Range("H137").Value = 1502.3823
$ws.Range("I137").Value = 1151.3889
$ws.Range("J137").Value = 1897.25
$ws.Range("K137").Value = 3454.1667
$ws.Range("L137").Value = 5691.75
$ws.Range("M137").Value = -904.1666999999998
$ws.Range("N137").Value = -10791.75
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1040.125
$ws.Range("I2").Value = 928.5833
$ws.Range("J2").Value = 1374.75
$ws.Range("K2").Value = 928.5833
$ws.Range("L2").Value = 1374.75
$ws.Range("M2").Value = -815.5833
$ws.Range("N2").Value = -1600.75
$ws.Range("H32").Value = 2719.79
$ws.Range("I32").Value = 2152.8823
$ws.Range("J32").Value = 5932.2666
$ws.Range("K32").Value = 2152.8823
$ws.Range("L32").Value = 5932.2666
$ws.Range("M32").Value = -1865.8823
$ws.Range("N32").Value = -6506.2666
$ws.Range("H102").Value = 16668296
$ws.Range("I102").Value = 20834620
$ws.Range("K102").Value = 20834620
$ws.Range("M102").Value = -20832998
$ws.Range("H116").Value = 1040.125
$ws.Range("I116").Value = 928.5833
$ws.Range("J116").Value = 1374.75
$ws.Range("K116").Value = 928.5833
$ws.Range("L116").Value = 1374.75
$ws.Range("M116").Value = 1365.4167
$ws.Range("N116").Value = -5962.75
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1040.125
$ws.Range("I3").Value = 928.5833
$ws.Range("J3").Value = 1374.75
$ws.Range("K3").Value = 928.5833
$ws.Range("L3").Value = 1374.75
$ws.Range("M3").Value = -814.5833
$ws.Range("N3").Value = -1602.75
$ws.Range("H82").Value = 19102.8
$ws.Range("J82").Value = 31000
$ws.Range("L82").Value = 31000
$ws.Range("N82").Value = -31766
$ws.Range("H85").Value = 19102.8
$ws.Range("J85").Value = 31000
$ws.Range("L85").Value = 31000
$ws.Range("N85").Value = -33652
$ws.Range("H99").Value = 90910344
$ws.Range("I99").Value = 111112280
$ws.Range("J99").Value = 1650
$ws.Range("K99").Value = 111112280
$ws.Range("L99").Value = 1650
$ws.Range("M99").Value = -111110782
$ws.Range("N99").Value = -4646
$ws.Range("H107").Value = 1534.4166
$ws.Range("I107").Value = 1088.8889
$ws.Range("J107").Value = 2871
$ws.Range("K107").Value = 1088.8889
$ws.Range("L107").Value = 2871
$ws.Range("M107").Value = 831.1111000000001
$ws.Range("N107").Value = -6711
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 10000
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 10000
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 10000
$ws.Range("N94").Value = -10902
$ws.Range("M94").ClearContents()
$ws.Range("H105").Value = 802.5
$ws.Range("I105").Value = 766.4286
$ws.Range("K105").Value = 766.4286
$ws.Range("M105").Value = 980.5714
$ws.Range("H111").Value = 40000
$ws.Range("J111").Value = 40000
$ws.Range("L111").Value = 40000
$ws.Range("N111").Value = -48180
$ws.Range("H112").Value = 36337.5
$ws.Range("J112").Value = 36337.5
$ws.Range("L112").Value = 36337.5
$ws.Range("N112").Value = -39291.5
$ws.Range("H132").Value = 1142.6765
$ws.Range("I132").Value = 809.0417
$ws.Range("K132").Value = 2427.1251
$ws.Range("M132").Value = 102.8748999999998
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 223.11765
$ws.Range("I14").Value = 223.11765
$ws.Range("K14").Value = 669.35295
$ws.Range("M14").Value = -496.35295
$ws.Range("H39").Value = 4083.6667
$ws.Range("J39").Value = 4140.4
$ws.Range("L39").Value = 12421.2
$ws.Range("N39").Value = -13009.2
$ws.Range("H129").Value = 18117060
$ws.Range("J129").Value = 4903198.5
$ws.Range("L129").Value = 14709595.5
$ws.Range("N129").Value = -14719595.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5146
$ws.Range("I80").Value = 6350
$ws.Range("J80").Value = 4458
$ws.Range("K80").Value = 6350
$ws.Range("L80").Value = 4458
$ws.Range("M80").Value = -5352
$ws.Range("N80").Value = -6454
$ws.Range("H83").Value = 5146
$ws.Range("I83").Value = 6350
$ws.Range("J83").Value = 4458
$ws.Range("K83").Value = 31750
$ws.Range("L83").Value = 22290
$ws.Range("M83").Value = -26758
$ws.Range("N83").Value = -32274
$ws.Range("H97").Value = 641.7646999999999
$ws.Range("I97").Value = 653.8889
$ws.Range("J97").Value = 628.125
$ws.Range("K97").Value = 653.8889
$ws.Range("L97").Value = 628.125
$ws.Range("M97").Value = -157.8889
$ws.Range("N97").Value = -1620.125
$ws.Range("H126").Value = 2138.75
$ws.Range("I126").Value = 1842.5
$ws.Range("J126").Value = 2435
$ws.Range("K126").Value = 5527.5
$ws.Range("L126").Value = 7305
$ws.Range("M126").Value = -3057.5
$ws.Range("N126").Value = -12245
$ws.Range("H132").Value = 3426.32
$ws.Range("I132").Value = 3625.6155
$ws.Range("J132").Value = 3210.4167
$ws.Range("K132").Value = 10876.8465
$ws.Range("L132").Value = 9631.250100000001
$ws.Range("M132").Value = -8346.8465
$ws.Range("N132").Value = -14691.2501
$ws.Range("H134").Value = 26987.889
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 26987.889
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 80963.667
$ws.Range("N134").Value = -86033.667
$ws.Range("M134").ClearContents()
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1237.8
$ws.Range("I100").Value = 1096.3334
$ws.Range("K100").Value = 1096.3334
$ws.Range("M100").Value = -555.3334
$ws.Range("H122").Value = 22729368
$ws.Range("I122").Value = 41668724
$ws.Range("J122").Value = 2141
$ws.Range("K122").Value = 125006172
$ws.Range("L122").Value = 6423
$ws.Range("M122").Value = -125003722
$ws.Range("N122").Value = -11323
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 38559.8
$ws.Range("J109").Value = 35614.25
$ws.Range("L109").Value = 35614.25
$ws.Range("N109").Value = -38388.25
$ws.Range("H113").Value = 529.5
$ws.Range("I113").Value = 364
$ws.Range("J113").Value = 695
$ws.Range("K113").Value = 1092
$ws.Range("L113").Value = 2085
$ws.Range("M113").Value = 1078
$ws.Range("N113").Value = -6425
$ws.Range("H132").Value = 1585
$ws.Range("I132").Value = 1370.8334
$ws.Range("J132").Value = 2227.5
$ws.Range("K132").Value = 4112.5002
$ws.Range("L132").Value = 6682.5
$ws.Range("M132").Value = -1582.5002
$ws.Range("N132").Value = -11742.5
$ws.Range("H133").Value = 33123
$ws.Range("J133").Value = 33123
$ws.Range("L133").Value = 33123
$ws.Range("N133").Value = -43243

Write-Output "Updated 184 cell value(s); cleared 2 obsolete cell(s)."
